$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cells for columns I (I0) and J (IF), formatted like the existing
# header row (bold font, thin border around, centered horizontally, top vertically)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$headerRange = $ws.Range("I1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous
$headerRange.Borders.Weight = 2            # xlThin

# Data values for column I (I0) and column J (IF)
$iValues = @{
    2 = 1; 3 = 1; 4 = 1; 5 = 1; 6 = 1; 7 = 1; 8 = 1; 9 = 1; 10 = 1;
    11 = 1; 12 = 1; 13 = 1; 14 = 1; 15 = 1; 16 = 1; 17 = 1; 18 = 1; 19 = 1; 20 = 1;
    21 = 1; 22 = 1; 23 = 1; 24 = 1; 25 = 1; 26 = 6; 27 = 1; 28 = 1; 29 = 1; 30 = 1
}

$jValues = @{
    2 = 4; 3 = 5; 4 = 4; 5 = 3; 6 = 3; 7 = 2; 8 = 6; 9 = 6; 10 = 6;
    11 = 7; 12 = 6; 13 = 6; 14 = 4; 15 = 6; 16 = 7; 17 = 3; 18 = 8; 19 = 6; 20 = 5;
    21 = 7; 22 = 5; 23 = 7; 24 = 5; 25 = 6; 26 = 9; 27 = 4; 28 = 3; 29 = 2; 30 = 2
}

for ($row = 2; $row -le 30; $row++) {
    $ws.Cells.Item($row, 9).Value = $iValues[$row]
    $ws.Cells.Item($row, 10).Value = $jValues[$row]
}

$wb.Save()
